# Weekly refresh of Fruta / hortaliza dataset (Membrillo - Vega Monumental Concepcion)
# The underlying data feed re-ordered the historical weekly observations;
# this script rewrites columns D (Fecha) and L:T (Calidad..Kg/unidad) for rows 2-27
# of the active sheet to match the refreshed dataset. Columns A:C, E:K are unchanged
# (same market/product/category/variety for every row in this subset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dData = @(
    @(44299),
    @(44299),
    @(44425),
    @(45041),
    @(45013),
    @(44698),
    @(45034),
    @(45014),
    @(45027),
    @(45029),
    @(44363),
    @(45079),
    @(44999),
    @(44999),
    @(44316),
    @(45092),
    @(45050),
    @(44776),
    @(44776),
    @(44358),
    @(45037),
    @(44307),
    @(44307),
    @(44272),
    @(44272),
    @(45076)
)

$ltData = @(
    @('Primera',100,10000,11000,10500,'$/caja 18 kilos granel','Región del Maule',583,18),
    @('Segunda',50,9000,9000,9000,'$/caja 18 kilos granel','Región del Maule',500,18),
    @('Primera',100,12000,13000,12500,'$/bandeja 18 kilos granel','Región de O''Higgins',694,18),
    @('Primera',100,11000,12000,11500,'$/bandeja 18 kilos granel','Región de O''Higgins',639,18),
    @('Primera',100,9000,10000,9500,'$/bandeja 18 kilos granel','Región de O''Higgins',528,18),
    @('Primera',50,10000,10000,10000,'$/caja 18 kilos granel','Región de O''Higgins',556,18),
    @('Primera',220,8500,9000,8727,'$/caja 18 kilos granel','Región de O''Higgins',485,18),
    @('Primera',100,9000,10000,9500,'$/bandeja 18 kilos granel','Región de O''Higgins',528,18),
    @('Primera',100,9000,10000,9500,'$/bandeja 18 kilos granel','Región de O''Higgins',528,18),
    @('Primera',100,9000,10000,9500,'$/bandeja 18 kilos granel','Región de O''Higgins',528,18),
    @('Primera',100,9000,10000,9500,'$/caja 15 kilos empedrada','Región de O''Higgins',633,15),
    @('Primera',270,11000,12000,11444,'$/caja 18 kilos granel','Región de O''Higgins',636,18),
    @('Primera',100,12000,12000,12000,'$/bandeja 18 kilos granel','Región de O''Higgins',667,18),
    @('Segunda',100,10000,10000,10000,'$/bandeja 18 kilos granel','Región de O''Higgins',556,18),
    @('Primera',100,9000,10000,9500,'$/caja 18 kilos granel','Región de O''Higgins',528,18),
    @('Primera',110,10000,11000,10455,'$/bandeja 18 kilos granel','Provincia de Curicó',581,18),
    @('Primera',140,11000,12000,11429,'$/caja 18 kilos empedrada','Región de O''Higgins',635,18),
    @('Primera',50,10000,10000,10000,'$/bandeja 18 kilos granel','Región de O''Higgins',556,18),
    @('Segunda',50,8000,8000,8000,'$/bandeja 18 kilos granel','Región de O''Higgins',444,18),
    @('Primera',100,11000,12000,11500,'$/caja 18 kilos granel','Región de O''Higgins',639,18),
    @('Primera',250,9000,9500,9200,'$/caja 18 kilos granel','Provincia de Curicó',511,18),
    @('Primera',50,10000,10000,10000,'$/bandeja 18 kilos granel','Región de O''Higgins',556,18),
    @('Segunda',50,8000,8000,8000,'$/bandeja 18 kilos granel','Región de O''Higgins',444,18),
    @('Primera',100,9000,10000,9500,'$/caja 15 kilos granel','Región de O''Higgins',633,15),
    @('Segunda',50,8000,8000,8000,'$/caja 15 kilos granel','Región de O''Higgins',533,15),
    @('Primera',150,10000,11000,10467,'$/caja 18 kilos granel','Provincia de Curicó',582,18)
)

$startRow = 2
for ($i = 0; $i -lt $dData.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value = $dData[$i][0]

    $row = $ltData[$i]
    $ws.Cells.Item($r, 12).Value = $row[0]   # L Calidad
    $ws.Cells.Item($r, 13).Value = $row[1]   # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[2]   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[3]   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[4]   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[5]   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $row[6]   # R Origen
    $ws.Cells.Item($r, 19).Value = $row[7]   # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[8]   # T Kg / unidad
}

